# Auto-generated edit script: update crypto price/volume table (cryptos.xlsx)
# Applies the hourly-refresh values from the GitHub Actions bot commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Coin/Link/Price/Volume cell in this sheet is stored as literal
# text (inline strings), even numeric-looking ones like '1.00' or
# '0.0000260'. Force columns B:E to Text format *before* writing so
# COM doesn't silently reinterpret them as numbers/percentages and
# normalize away significant trailing zeros / thousands-dot formatting.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '68.333.54'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '3.744.46'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '593.10'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').Value = '166.18'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '3.742.85'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '0.0000260'
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('D14').Value = '36.15'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '4.372.65'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '3.751.49'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '68.283.91'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').Value = '17.85'
$ws.Range('E18').Value = '  -3.04%  '
$ws.Range('D19').Value = '6.99'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = '10.65'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('D22').Value = '464.48'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '0.695'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '83.94'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000147'
$ws.Range('E25').Value = '  +4.06%  '
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('D27').Value = '11.87'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').Value = '  -2.38%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '3.892.23'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '2.76'
$ws.Range('E31').Value = '  -4.41%  '
$ws.Range('D32').Value = '7.29'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '29.86'
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').Value = '9.17'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('D37').Value = '3.700.50'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('D39').Value = '3.43'
$ws.Range('E39').Value = '  -4.27%  '
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = '5.78'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '0.301'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '46.71'
$ws.Range('E46').Value = '  +4.04%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Value = '42.87'
$ws.Range('E48').Value = '  +10.08%  '
$ws.Range('D49').Value = '8.46'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').Value = '388.46'
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('D51').Value = '144.02'
$ws.Range('E51').Value = '  +0.10%  '
